$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (13) to the routine log: copy the date cell's formatting
# from the row above (A12) so the new date cell uses the same style,
# then set the actual values for the new entry.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A13").Value = 43819
$ws.Range("B13").Value = "Full stack management"
